$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds "xyz_path" values for rows 2-11 (H2_0 .. H2_9).
# Replace the old relative path segment with the new one.
$oldSeg = "tests\test_simulation\test_optimizers"
$newSeg = "gym\simulation\optimizers"

for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $oldVal = [string]$cell.Value()
    $newVal = $oldVal.Replace($oldSeg, $newSeg)
    $cell.Value = $newVal
}
